# Add 2022-Q1 sheet (feat: add 2022-Q1 data)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

# Writes $value into $range forcing a TEXT (string) cell, with no left-over
# numeric-format styling (mirrors the original workbook where such cells
# carry no explicit style index).
function Set-TextValue($range, $blank, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $blank.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats - restores the default (no) style
}

# Writes a numeric $value into $range with no left-over styling.
function Set-NumberValue($range, $blank, $value) {
    $range.Value = $value
    $blank.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" and before "总计"
# ---------------------------------------------------------------------------

$wsQ4 = $wb.Worksheets.Item("2021-Q4")

$wsNew = $wb.Worksheets.Add($null, $wsQ4)
$wsNew.Name = "2022-Q1"
$wsNew.Outline.SummaryRow = 1
$wsNew.Outline.SummaryColumn = 1

# NOTE: worksheet references resolve dynamically by position, so "总计" must
# be (re-)looked-up *after* the new sheet has been inserted - otherwise the
# old reference would silently now point at the newly inserted sheet.
$wsTotal = $wb.Worksheets.Item("总计")

# Unformatted cell used as the "blank" format source for PasteSpecial calls.
$blankNew = $wsNew.Range("Z1000")

# ---------------------------------------------------------------------------
# 2. Header row (B1:H1) - bold / bordered / centered, same style as sheet "2021-Q4"
# ---------------------------------------------------------------------------

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsNew.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
}

# Copy the header formatting (bold, border, centered) from the "2021-Q4" sheet.
$wsQ4.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------------

# Row 2
Set-NumberValue $wsNew.Range("A2") $blankNew 0
$wsQ4.Range("A2").Copy()
$wsNew.Range("A2").PasteSpecial(-4122)

Set-TextValue $wsNew.Range("B2") $blankNew "090010"
Set-TextValue $wsNew.Range("C2") $blankNew "大成中证红利指数A"
Set-TextValue $wsNew.Range("D2") $blankNew "34.51"
Set-TextValue $wsNew.Range("E2") $blankNew "93.73"
Set-TextValue $wsNew.Range("F2") $blankNew "1.32"
Set-TextValue $wsNew.Range("G2") $blankNew "0.4555"
Set-NumberValue $wsNew.Range("H2") $blankNew 10

# Row 3
Set-NumberValue $wsNew.Range("A3") $blankNew 1
$wsQ4.Range("A2").Copy()
$wsNew.Range("A3").PasteSpecial(-4122)

Set-TextValue $wsNew.Range("B3") $blankNew "007801"
Set-TextValue $wsNew.Range("C3") $blankNew "大成中证红利指数C"
Set-TextValue $wsNew.Range("D3") $blankNew "3.87"
Set-TextValue $wsNew.Range("E3") $blankNew "93.73"
Set-TextValue $wsNew.Range("F3") $blankNew "1.32"
Set-TextValue $wsNew.Range("G3") $blankNew "0.0511"
Set-NumberValue $wsNew.Range("H3") $blankNew 10

# ---------------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: insert a new row for 2022-Q1 above the
#    existing 2021-Q4 row.
# ---------------------------------------------------------------------------

$wsTotal.Rows(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$blankTotal = $wsTotal.Range("Z1000")

Set-NumberValue $wsTotal.Range("A2") $blankTotal 0
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

Set-TextValue $wsTotal.Range("B2") $blankTotal "2022-Q1"
Set-NumberValue $wsTotal.Range("C2") $blankTotal 2
Set-NumberValue $wsTotal.Range("D2") $blankTotal 0.51

# The previously existing row (2021-Q4) shifted down to row 3; restore its
# index value ("1") which the row-shift does not update automatically.
$wsTotal.Range("A3").Value = 1
